$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing (pre-edit) data rows 2..8, columns A..D, into scalars
# before overwriting anything, since the new order interleaves old rows.
$A2 = $ws.Cells.Item(2,1).Value2; $B2 = $ws.Cells.Item(2,2).Value2; $C2 = $ws.Cells.Item(2,3).Value2; $D2 = $ws.Cells.Item(2,4).Value2
$A3 = $ws.Cells.Item(3,1).Value2; $B3 = $ws.Cells.Item(3,2).Value2; $C3 = $ws.Cells.Item(3,3).Value2; $D3 = $ws.Cells.Item(3,4).Value2
$A4 = $ws.Cells.Item(4,1).Value2; $B4 = $ws.Cells.Item(4,2).Value2; $C4 = $ws.Cells.Item(4,3).Value2; $D4 = $ws.Cells.Item(4,4).Value2
$A5 = $ws.Cells.Item(5,1).Value2; $B5 = $ws.Cells.Item(5,2).Value2; $C5 = $ws.Cells.Item(5,3).Value2; $D5 = $ws.Cells.Item(5,4).Value2
$A6 = $ws.Cells.Item(6,1).Value2; $B6 = $ws.Cells.Item(6,2).Value2; $C6 = $ws.Cells.Item(6,3).Value2; $D6 = $ws.Cells.Item(6,4).Value2
$A7 = $ws.Cells.Item(7,1).Value2; $B7 = $ws.Cells.Item(7,2).Value2; $C7 = $ws.Cells.Item(7,3).Value2; $D7 = $ws.Cells.Item(7,4).Value2
$A8 = $ws.Cells.Item(8,1).Value2; $B8 = $ws.Cells.Item(8,2).Value2; $C8 = $ws.Cells.Item(8,3).Value2; $D8 = $ws.Cells.Item(8,4).Value2

# The rows are being re-sorted by ascending time (column A):
#   new row 2 <- old row 4
#   new row 3 <- old row 7
#   new row 4 <- old row 5
#   new row 5 <- old row 3
#   new row 6 <- old row 8
#   new row 7 <- old row 6
#   new row 8 <- old row 2
$ws.Cells.Item(2,1).Value = $A4; $ws.Cells.Item(2,2).Value = $B4; $ws.Cells.Item(2,3).Value = $C4; $ws.Cells.Item(2,4).Value = $D4
$ws.Cells.Item(3,1).Value = $A7; $ws.Cells.Item(3,2).Value = $B7; $ws.Cells.Item(3,3).Value = $C7; $ws.Cells.Item(3,4).Value = $D7
$ws.Cells.Item(4,1).Value = $A5; $ws.Cells.Item(4,2).Value = $B5; $ws.Cells.Item(4,3).Value = $C5; $ws.Cells.Item(4,4).Value = $D5
$ws.Cells.Item(5,1).Value = $A3; $ws.Cells.Item(5,2).Value = $B3; $ws.Cells.Item(5,3).Value = $C3; $ws.Cells.Item(5,4).Value = $D3
$ws.Cells.Item(6,1).Value = $A8; $ws.Cells.Item(6,2).Value = $B8; $ws.Cells.Item(6,3).Value = $C8; $ws.Cells.Item(6,4).Value = $D8
$ws.Cells.Item(7,1).Value = $A6; $ws.Cells.Item(7,2).Value = $B6; $ws.Cells.Item(7,3).Value = $C6; $ws.Cells.Item(7,4).Value = $D6
$ws.Cells.Item(8,1).Value = $A2; $ws.Cells.Item(8,2).Value = $B2; $ws.Cells.Item(8,3).Value = $C2; $ws.Cells.Item(8,4).Value = $D2
